$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToDO")

# Add two new "To Do" entries into rows 12 and 13 (mock-up review task)
$ws.Range("A12").Value = "Nazicht mock-ups"
$ws.Range("B12").Value = "15 minuten"
$ws.Range("C12").Value = "30 minuten"
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = 41363
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = "Steven V"
$ws.Range("G12").Value = "Solved"
$ws.Range("H12").Value = "APP"

$ws.Range("A13").Value = "Nazicht mock-ups"
$ws.Range("B13").Value = "15 minuten"
$ws.Range("C13").Value = "15 minuten"
$ws.Range("D11").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = 41364
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = "Steven V"
$ws.Range("G13").Value = "Solved"
$ws.Range("H13").Value = "APP"

# Move the active selection to A14, matching the saved cursor position
$ws.Range("A14").Select()

$excel.Calculate()
